# Updated value set parser to handle new MAT valueset format.
#
# 1) The "Disclaimer" sheet is an empty placeholder tab that is removed
#    entirely.
# 2) "White_List" becomes the first (and active/selected) sheet.
# 3) The old "QDM Category" column (column E) in White_List is dropped -
#    deleting it shifts Code System / Code System Version / Code /
#    Descriptor / Concept one column to the left.
# 4) The worksheet's remembered selection moves from A6 to G32.

$wb = $excel.ActiveWorkbook

# Remove the empty "Disclaimer" sheet.
$wb.Worksheets.Item("Disclaimer").Delete()

$ws = $wb.Worksheets.Item("White_List")

# Make White_List the active/selected sheet (it is now first in tab order).
$ws.Activate()

# Drop the obsolete "QDM Category" column; everything to the right shifts left.
$ws.Columns("E").Delete()

# Restore the saved selection on the sheet.
$ws.Range("G32").Select()
